$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 11 (pushes the old GBS row from 11 -> 12),
# making room for the new "HUS / mild" row.
$ws.Rows(11).Insert()

# --- Header row: add "unit" column header ---
$ws.Range("D1").Value = "unit"

# --- GI (per day) ---
$ws.Range("D2").Value = "per day"
$ws.Range("D3").Value = "per day"

# --- flulike (per day) ---
$ws.Range("D4").Value = "per day"
$ws.Range("D5").Value = "per day"

# --- IBS severe/mild: drop the /365 division, now annual values ---
$ws.Range("C6").Formula = "=964"
$ws.Range("D6").Value = "per year"

$ws.Range("C7").Formula = "=344"
$ws.Range("D7").Value = "per year"

# --- ReA severe/mild ---
$ws.Range("C8").Formula = "=1166"
$ws.Range("D8").Value = "per year"

$ws.Range("C9").Formula = "=605"
$ws.Range("D9").Value = "per year"

# --- HUS severe (existing row 10) ---
$ws.Range("C10").Formula = "=1620"
$ws.Range("D10").Value = "per year"

# --- HUS mild (new row 11) ---
$ws.Range("A11").Value = "HUS"
$ws.Range("B11").Value = "mild"
$ws.Range("C11").Formula = "=901"
$ws.Range("D11").Value = "per year"

# --- GBS severe (was row 11, now shifted to row 12) ---
$ws.Range("A12").Value = "GBS"
$ws.Range("B12").Value = "severe"
$ws.Range("C12").Formula = "=1371"
$ws.Range("D12").Value = "per year"

# --- GBS mild (new row 13) ---
$ws.Range("A13").Value = "GBS"
$ws.Range("B13").Value = "mild"
$ws.Range("C13").Formula = "=762"
$ws.Range("D13").Value = "per year"

# Match the author's final cursor position left over from editing.
[void]$ws.Range("F6").Select()
